$wb = $excel.ActiveWorkbook

# --- Sheet "compare_models": update TT (Sec) column (I) for several rows ---
$ws1 = $wb.Worksheets.Item("compare_models")

$ws1.Range("I2").Value  = 0.058
$ws1.Range("I3").Value  = 0.074
$ws1.Range("I4").Value  = 0.032
$ws1.Range("I5").Value  = 0.092
$ws1.Range("I6").Value  = 0.046
$ws1.Range("I7").Value  = 1.028
$ws1.Range("I9").Value  = 0.018
$ws1.Range("I10").Value = 0.026
$ws1.Range("I11").Value = 0.534
$ws1.Range("I12").Value = 0.02
$ws1.Range("I14").Value = 0.016
$ws1.Range("I16").Value = 0.018
$ws1.Range("I17").Value = 0.016
$ws1.Range("I18").Value = 0.018

# --- Sheet "pred_final": update metrics row 2 (C:H) ---
$ws2 = $wb.Worksheets.Item("pred_final")

$ws2.Range("C2").Value = 1.5421
$ws2.Range("D2").Value = 3.8317
$ws2.Range("E2").Value = 1.9575
$ws2.Range("F2").Value = 0.9943
$ws2.Range("G2").Value = 0.0384
$ws2.Range("H2").Value = 0.0278
